# feat: add 2022-Q4 data
#
# 1. Insert a new sheet "2022-Q4" right after "总计" (position 2), pushing
#    the existing "2022-Q3" / "2022-Q1" / "2021-Q4" sheets one slot to the
#    right.
# 2. Populate the new "2022-Q4" sheet with the Q4 fund-holding table.
# 3. Insert a new row into "总计" summarising the 2022-Q4 quarter and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new worksheet right after "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Populate "2022-Q4" - same layout/styles as the "2022-Q3" sheet
# ---------------------------------------------------------------------

# Header row (B1:H1) - reuse the bordered/bold header style already used
# on sheet "2022-Q3".
$q3.Cells.Item(1, 2).Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Row-index column (A2:A7) - reuse the bordered/bold index style.
$q3.Cells.Item(2, 1).Copy()
$q4.Range("A2:A7").PasteSpecial(-4122)

# The D:G columns store numeric-looking figures as plain text (exactly
# like the other quarters' sheets), so force a Text number format before
# assigning the values to stop them turning into real numbers.
$q4.Range("B2:G7").NumberFormat = "@"

$rows = @(
    @("005585", "银河文体娱乐主题灵活配置混合A", "3.15", "88.90", "9.19", "0.2895", 1),
    @("013890", "国泰睿毅三年持有期混合A",       "4.82", "90.04", "3.80", "0.1832", 8),
    @("001628", "招商体育文化休闲股票A",         "2.33", "93.03", "5.47", "0.1275", 2),
    @("015667", "银河文体娱乐主题灵活配置混合C", "0.38", "88.90", "9.19", "0.0349", 1),
    @("013891", "国泰睿毅三年持有期混合C",       "0.45", "90.04", "3.80", "0.0171", 8),
    @("015395", "招商体育文化休闲股票C",         "0.29", "93.03", "5.47", "0.0159", 2)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $q4.Cells.Item($excelRow, 1).Value = $r
    $q4.Cells.Item($excelRow, 2).Value = $row[0]
    $q4.Cells.Item($excelRow, 3).Value = $row[1]
    $q4.Cells.Item($excelRow, 4).Value = $row[2]
    $q4.Cells.Item($excelRow, 5).Value = $row[3]
    $q4.Cells.Item($excelRow, 6).Value = $row[4]
    $q4.Cells.Item($excelRow, 7).Value = $row[5]
    $q4.Cells.Item($excelRow, 8).Value = $row[6]
}

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Insert the new summary row into "总计"
# ---------------------------------------------------------------------

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.67

# Restore the bordered/bold style on A2 (Insert copies row-1's header
# style, not the A-column data style).
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

# Renumber the pushed-down rows' index column.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3

$total.Range("A1").Select()
